# ranking.xlsx — "Add files via upload"
#
# Appends 20 new ranking-table rows (216-235) under the existing A:B
# columns (Name / Points). One of the rows references a brand-new
# competitor, "שלו דיין", who is not yet present in the shared-string
# table; simply writing her name as a cell value will make Excel add
# a new shared string for her automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: name, Column B: points — in row order 216..235
$names = @(
    "הגר אגמון",
    "תאיו ורד",
    "אורי שטרנברג",
    "יובל סטרוזר",
    "איתי הראל",
    "יובל סטרוזר",
    "איתי הראל",
    "עדן ורד מרי",
    "אביב ואסקז",
    "יער אלביר",
    "איתי הראל",
    "הילה שולויס",
    "איתי בסטקר",
    "ליאם דיין ",
    "לינוי קוסטיקה",
    "שלו דיין",
    "ליהי בראל",
    "יולי יערי תליו",
    "לינוי קוסטיקה",
    "איתי הראל"
)

$points = @(1, 1, 1, 1, 1, 6, 6, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 6, 6)

$startRow = 216
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $points[$i]
}

# Mirror the author's final view state (scrolled down, cell B224 selected).
$ws.Range("A213").Select()
$ws.Range("B224").Select()
